$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date values for rows 2-5 from 45233 to 45243
$ws.Range("C2").Value = 45243
$ws.Range("C3").Value = 45243
$ws.Range("C4").Value = 45243
$ws.Range("C5").Value = 45243
